$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 10999.714
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 10999.714
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 10999.714
$ws.Range("N32").Value = -11651.714

# Row 40
$ws.Range("H40").Value = 8153.276
$ws.Range("I40").Value = 4099.4
$ws.Range("J40").Value = 8997.833000000001
$ws.Range("K40").Value = 4099.4
$ws.Range("L40").Value = 8997.833000000001
$ws.Range("M40").Value = -3924.4
$ws.Range("N40").Value = -9347.833000000001

# Row 64
$ws.Range("H64").Value = 50002500
$ws.Range("I64").Value = 100000000
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 100000000
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -99999752
$ws.Range("N64").Value = -5496

# Row 67
$ws.Range("H67").Value = 50002500
$ws.Range("I67").Value = 100000000
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 100000000
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -99999142
$ws.Range("N67").Value = -6716

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 8000
$ws.Range("I45").Value = 8000
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 8000
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -7623

# Row 97
$ws.Range("H97").Value = 539.44446
$ws.Range("I97").Value = 549.1429000000001
$ws.Range("J97").Value = 505.5
$ws.Range("K97").Value = 549.1429000000001
$ws.Range("L97").Value = 505.5
$ws.Range("M97").Value = -53.14290000000005

# Row 102
$ws.Range("H102").Value = 3959.8
$ws.Range("I102").Value = 4199.75
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 4199.75
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -2577.75
$ws.Range("N102").Value = -6244

# Row 122
$ws.Range("H122").Value = 4200

# Row 132
$ws.Range("H132").Value = 1447.5
$ws.Range("I132").Value = 1447.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4342.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1812.5
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 4004.5
$ws.Range("I99").Value = 3010
$ws.Range("J99").Value = 4999
$ws.Range("K99").Value = 3010
$ws.Range("L99").Value = 4999
$ws.Range("M99").Value = -1512
$ws.Range("N99").Value = -7995

# Row 105
$ws.Range("H105").Value = 4816.1665
$ws.Range("I105").Value = 5333
$ws.Range("J105").Value = 4299.3335
$ws.Range("K105").Value = 5333
$ws.Range("L105").Value = 4299.3335
$ws.Range("M105").Value = -3586
$ws.Range("N105").Value = -7793.3335

# Row 134
$ws.Range("H134").Value = 3399.7144
$ws.Range("I134").Value = 3399.7144
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 10199.1432
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7664.143199999999

$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 373.875
$ws.Range("I10").Value = 487.33334
$ws.Range("J10").Value = 33.5
$ws.Range("K10").Value = 487.33334
$ws.Range("L10").Value = 33.5
$ws.Range("M10").Value = -348.33334
$ws.Range("N10").Value = -311.5

# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

# Row 31
$ws.Range("H31").Value = 2024.0952
$ws.Range("I31").Value = 710.9167
$ws.Range("J31").Value = 3775
$ws.Range("K31").Value = 710.9167
$ws.Range("L31").Value = 3775
$ws.Range("M31").Value = -415.9167

# Row 34
$ws.Range("H34").Value = 2024.0952
$ws.Range("I34").Value = 710.9167
$ws.Range("J34").Value = 3775
$ws.Range("K34").Value = 710.9167
$ws.Range("L34").Value = 3775
$ws.Range("M34").Value = -508.9167

# Row 50
$ws.Range("H50").Value = 23222.223
$ws.Range("I50").Value = 26666.666
$ws.Range("J50").Value = 21500
$ws.Range("K50").Value = 26666.666
$ws.Range("L50").Value = 21500
$ws.Range("M50").Value = -26041.666
$ws.Range("N50").Value = -22750

# Row 51
$ws.Range("H51").Value = 25000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 25000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 25000
$ws.Range("N51").Value = -26472
$ws.Range("M51").ClearContents()

# Row 60
$ws.Range("H60").Value = 13170.588
$ws.Range("I60").Value = 8990.909
$ws.Range("J60").Value = 20833.334
$ws.Range("K60").Value = 8990.909
$ws.Range("L60").Value = 20833.334
$ws.Range("M60").Value = -8479.909
$ws.Range("N60").Value = -21855.334

# Row 61
$ws.Range("H61").Value = 25000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 25000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 25000
$ws.Range("N61").Value = -25696
$ws.Range("M61").ClearContents()

# Row 99
$ws.Range("H99").Value = 1012
$ws.Range("I99").Value = 1012
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1012
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 486

# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

# Row 126
$ws.Range("H126").Value = 1012
$ws.Range("I126").Value = 1012
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3036
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -566

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 438.85715
$ws.Range("I12").Value = 1000.3333
$ws.Range("J12").Value = 17.75
$ws.Range("K12").Value = 3000.9999
$ws.Range("L12").Value = 53.25
$ws.Range("M12").Value = -2827.9999

# Row 25
$ws.Range("H25").Value = 198
$ws.Range("I25").Value = 198
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 594
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -425

# Row 30
$ws.Range("H30").Value = 198
$ws.Range("I30").Value = 198
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 594
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -492

# Row 64
$ws.Range("H64").Value = 1000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 1000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3000
$ws.Range("N64").Value = -3540
$ws.Range("M64").ClearContents()

# Row 67
$ws.Range("H67").Value = 1000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 1000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3000
$ws.Range("N67").Value = -4872
$ws.Range("M67").ClearContents()

# Row 75
$ws.Range("H75").Value = 490
$ws.Range("I75").Value = 490
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 1470
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -472

# Row 78
$ws.Range("H78").Value = 490
$ws.Range("I78").Value = 490
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 4410
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = 582

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3619.8
$ws.Range("I80").Value = 3450
$ws.Range("J80").Value = 3733
$ws.Range("K80").Value = 3450
$ws.Range("L80").Value = 3733
$ws.Range("M80").Value = -2452

# Row 83
$ws.Range("H83").Value = 3619.8
$ws.Range("I83").Value = 3450
$ws.Range("J83").Value = 3733
$ws.Range("K83").Value = 17250
$ws.Range("L83").Value = 18665
$ws.Range("M83").Value = -12258

# Row 86
$ws.Range("H86").Value = 45000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 45000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 45000
$ws.Range("N86").Value = -47372

# Row 89
$ws.Range("H89").Value = 45000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 45000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 135000
$ws.Range("N89").Value = -146856

# Row 102
$ws.Range("H102").Value = 3311.4546
$ws.Range("I102").Value = 3565.7
$ws.Range("J102").Value = 769
$ws.Range("K102").Value = 3565.7
$ws.Range("L102").Value = 769
$ws.Range("M102").Value = -1943.7

# Row 113
$ws.Range("H113").Value = 1825
$ws.Range("I113").Value = 1450
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 1450
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = 720
$ws.Range("N113").Value = -6540

# Row 126
$ws.Range("H126").Value = 3258.3125
$ws.Range("I126").Value = 1654.125
$ws.Range("J126").Value = 4862.5
$ws.Range("K126").Value = 4962.375
$ws.Range("L126").Value = 14587.5
$ws.Range("M126").Value = -2492.375

$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 11601
$ws.Range("I20").Value = 8005
$ws.Range("J20").Value = 12500
$ws.Range("K20").Value = 8005
$ws.Range("L20").Value = 12500
$ws.Range("M20").Value = -7779
$ws.Range("N20").Value = -12952

# Row 22
$ws.Range("H22").Value = 2800.1667
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 3260.2
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 3260.2
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -3850.2

# Row 24
$ws.Range("H24").Value = 15668.667
$ws.Range("I24").Value = 10006
$ws.Range("J24").Value = 18500
$ws.Range("K24").Value = 10006
$ws.Range("L24").Value = 18500
$ws.Range("M24").Value = -9663
$ws.Range("N24").Value = -19186

# Row 27
$ws.Range("H27").Value = 2800.1667
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 3260.2
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 3260.2
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -3474.2

# Row 40
$ws.Range("H40").Value = 6575.7036
$ws.Range("I40").Value = 7445.1816
$ws.Range("J40").Value = 2750
$ws.Range("K40").Value = 7445.1816
$ws.Range("L40").Value = 2750
$ws.Range("M40").Value = -7309.1816
$ws.Range("N40").Value = -3022

# Row 46
$ws.Range("H46").Value = 3790.182
$ws.Range("I46").Value = 2478.2
$ws.Range("J46").Value = 4883.5
$ws.Range("K46").Value = 2478.2
$ws.Range("L46").Value = 4883.5
$ws.Range("M46").Value = -2290.2
$ws.Range("N46").Value = -5259.5

# Row 61
$ws.Range("H61").Value = 4933.3335
$ws.Range("I61").Value = 5400
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 5400
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -5198
$ws.Range("N61").Value = -4404

# Row 113
$ws.Range("H113").Value = 4933.3335
$ws.Range("I113").Value = 5400
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 5400
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -3230
$ws.Range("N113").Value = -8340

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4744.6
$ws.Range("I62").Value = 4900
$ws.Range("J62").Value = 4705.75
$ws.Range("K62").Value = 4900
$ws.Range("L62").Value = 4705.75
$ws.Range("M62").Value = -4276
$ws.Range("N62").Value = -5953.75

# Row 65
$ws.Range("H65").Value = 4744.6
$ws.Range("I65").Value = 4900
$ws.Range("J65").Value = 4705.75
$ws.Range("K65").Value = 24500
$ws.Range("L65").Value = 23528.75
$ws.Range("M65").Value = -21380
$ws.Range("N65").Value = -29768.75

# Row 100
$ws.Range("H100").Value = 599
$ws.Range("I100").Value = 499
$ws.Range("J100").Value = 699
$ws.Range("K100").Value = 998
$ws.Range("L100").Value = 1398
$ws.Range("M100").Value = -457
$ws.Range("N100").Value = -2480

# Row 132
$ws.Range("H132").Value = 2506.1538
$ws.Range("I132").Value = 2480.7
$ws.Range("J132").Value = 2591
$ws.Range("K132").Value = 7442.099999999999
$ws.Range("L132").Value = 7773
$ws.Range("M132").Value = -4912.099999999999

# Row 136
$ws.Range("H136").Value = 3729.875
$ws.Range("I136").Value = 3141.75
$ws.Range("J136").Value = 5494.25
$ws.Range("K136").Value = 9425.25
$ws.Range("L136").Value = 16482.75
$ws.Range("M136").Value = -6875.25
